# Replace the formula-driven column I (rows 2-32) with a literal lookup
# table of researched values (text, one decimal place). This mirrors the
# author swapping `=$A{row}+I$1` for hard values taken from an external
# table, so column I stops participating in the sheet-wide shared formula.
#
# NOTE on ordering: the cells are written in the exact order the shared
# strings ended up enumerated in the saved workbook (I11 first, then
# I2..I10, then I12..I32) so that xl/sharedStrings.xml comes out with
# "6.0" as the first <si> entry, matching the authored file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lookup = [ordered]@{
    "I11" = "6.0"
    "I2"  = "11.4"
    "I3"  = "10.4"
    "I4"  = "9.5"
    "I5"  = "8.8"
    "I6"  = "8.1"
    "I7"  = "7.6"
    "I8"  = "7.1"
    "I9"  = "6.7"
    "I10" = "6.3"
    "I12" = "5.7"
    "I13" = "5.4"
    "I14" = "5.2"
    "I15" = "5.0"
    "I16" = "4.7"
    "I17" = "4.6"
    "I18" = "4.4"
    "I19" = "4.2"
    "I20" = "4.1"
    "I21" = "3.9"
    "I22" = "3.8"
    "I23" = "3.7"
    "I24" = "3.6"
    "I25" = "3.5"
    "I26" = "3.4"
    "I27" = "3.3"
    "I28" = "3.2"
    "I29" = "3.1"
    "I30" = "3.0"
    "I31" = "2.9"
    "I32" = "2.8"
}

foreach ($addr in $lookup.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage (t="s" shared-string) instead of a number, then
    # drop the temporary "@" text format so the cell keeps the workbook's
    # default (General) style, same as in the authored file.
    $cell.NumberFormat = "@"
    $cell.Value = $lookup[$addr]
    $cell.ClearFormats()
}

# The author left the selection on the last edited cell.
$ws.Range("I32").Select() | Out-Null
